# Generate Report for Handback
# Update the generated-report timestamps on the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 74132d39-6caf-4ef3-9913-a0329f4edc59.md
$wsOverview.Range("G3").Value = "2016-08-25 14:46:45"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the same file
$wsZhCn.Range("H3").Value = "2016-08-25 14:46:41"
$wsZhCn.Range("K3").Value = "2016-08-25 14:47:01"

# de-de sheet: "Correspond Handback DateTime" for the same file
$wsDeDe.Range("K3").Value = "2016-08-25 14:47:19"
